$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 1044 (shifts existing rows 1044.. down by one)
$ws.Rows.Item(1044).Insert()

# Populate the newly inserted row with this week's data
$ws.Range("A1044").Value = 6
$ws.Range("B1044").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1044").Value = "Metropolitana"
$ws.Range("D1044").Value = 45041
$ws.Range("E1044").Value = 13
$ws.Range("F1044").Value = 100112003
$ws.Range("G1044").Value = "Ajo"
$ws.Range("H1044").Value = "Chino"
$ws.Range("I1044").Value = "Primera"
$ws.Range("J1044").Value = 1400
$ws.Range("K1044").Value = 13500
$ws.Range("L1044").Value = 14000
$ws.Range("M1044").Value = 13679
$ws.Range("N1044").Value = "$/caja 10 kilos"
$ws.Range("O1044").Value = "China"
$ws.Range("P1044").Value = 1368
$ws.Range("Q1044").Value = 10
$ws.Range("R1044").Value = "Hortaliza"
